$d = $word.ActiveDocument

# 1. "Ator principal:" Contractor -> Contratante ou Ambos
$d.Content.Find.Execute("Contractor", $true, $false, $false, $false, $false, $true, 1, $false, "Contratante ou Ambos", 2) | Out-Null

# 2. Pre-condicoes: "- O ator estar “logado” no sistema." -> "- O ator estar autenticado" + " no sistema." (two runs)
$d.Content.Find.Execute("- O ator estar “logado” no sistema.", $true, $false, $false, $false, $false, $true, 1, $false, "- O ator estar autenticado no sistema.", 2) | Out-Null

# 3. Pos-condicoes paragraph: "Parceiro(s) adicionado(s) a lista de parcerias do ator." -> "Parceiro(a) adicionado(a) a lista de parcerias do ator."
$d.Content.Find.Execute("arceiro(s)", $true, $false, $false, $false, $false, $true, 1, $false, "arceiro(a)", 2) | Out-Null
$d.Content.Find.Execute("adicionado(s) a lista de parcerias do ator", $true, $false, $false, $false, $false, $true, 1, $false, "adicionado(a) a lista de parcerias do ator", 2) | Out-Null

# 4. Remove the "Fluxo principal" header row and its single content row (rows 7 and 8 of table 1)
$t = $d.Tables.Item(1)
$t.Rows.Item(7).Delete()
$t.Rows.Item(7).Delete()

# 5. Rename "Fluxo alternativo 1" header -> "Fluxo Principal"
$d.Content.Find.Execute("Fluxo alternativo 1", $true, $false, $false, $false, $false, $true, 1, $false, "Fluxo Principal", 2) | Out-Null

# 6. "...visualizar o perfil do mesmo." -> "...visualizar o perfil do possível novo parceiro."
$d.Content.Find.Execute("visualizar o perfil do mesmo.", $true, $false, $false, $false, $false, $true, 1, $false, "visualizar o perfil do possível novo parceiro.", 2) | Out-Null

# 7. "cujo o cliente escolheu para visualizar" -> "cujo o ator selecionou para visualizar"
$d.Content.Find.Execute("cujo o cliente escolheu para visualizar", $true, $false, $false, $false, $false, $true, 1, $false, "cujo o ator selecionou para visualizar", 2) | Out-Null

# 8. "clique no botão “Requisitar parceria” para adicionar este à sua lista de parceiros." -> "clique no botão “Adicionar Parceiro” para adicioná-lo  à sua lista de parceiros."
$d.Content.Find.Execute("clique no botão “Requisitar parceria” para adicionar este à sua lista de parceiros.", $true, $false, $false, $false, $false, $true, 1, $false, "clique no botão “Adicionar Parceiro” para adicioná-lo  à sua lista de parceiros.", 2) | Out-Null

# 9. "para aguardar a aprovação ou não do contratante requisitado e apresenta uma mensagem de sucesso na requisição." -> "...do parceiro requisitado..." (only within the remaining Fluxo Principal / Sub-fluxo1 section, last occurrence)
$d.Content.Find.Execute("não do contratante requisitado e apresenta", $true, $false, $false, $false, $false, $true, 1, $false, "não do parceiro requisitado e apresenta", 2) | Out-Null

# 10. Re-insert the "_GoBack" bookmark right after "parceiro" in that same sentence (collapsed bookmark),
#     mirroring where it moved to in the diff.
$fr = $d.Content
$fr.Find.Execute("parceiro requisitado", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPos = $fr.Start + 8
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
